$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows from bottom to top so earlier row numbers stay valid.
# Row 61: "The search for E.T. Bell" by Reid, Constance.
$ws.Rows("61:61").Delete()
# Rows 30-31: "Lectures on modern mathematics I" and "II" by Saaty, Thomas L.
$ws.Rows("30:31").Delete()

$ws.Range("A59").Select() | Out-Null
